$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.940.12'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.554.09'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.88%  '
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.79'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.65%  '
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0859'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").Value = '1.775.38'
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("D13").Value = '1.556.12'
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("E15").Value = '  +2.39%  '
$ws.Range("D16").Value = '26.929.58'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.68'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '216.98'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.23'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.10'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.96%  '
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0469'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.15%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").Value = '1.439.23'
$ws.Range("E33").Value = '  +5.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.03'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.81%  '
$ws.Range("E35").Value = '  +3.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.962'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.811'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.987'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("E44").Value = '  +3.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.84'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.94%  '
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("D47").Value = '1.689.63'
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.23'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0524'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("D50").Value = '0.0₇0997'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0957'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.64%  '
